{"js": "// The edit rewrites the whole body: the 2-paragraph \"void setup(){...} void\n// loop(){...}\" sketch becomes a longer sketch that declares redLED/greenLED/\n// blinkTimes, calls a new blink() helper from loop(), and defines blink()\n// itself. Because so many runs/paragraphs are inserted (and the diff shows\n// exact <w:proofErr> spell-check bracketing around every camelCase\n// identifier), the most faithful way to reproduce the target OOXML is to\n// build that OOXML and hand it to Body.insertOoxml with Replace semantics.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst ooxml = `<?xml version=\"1.0\" standalone=\"yes\"?><?mso-application progid=\"Word.Document\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>redLED</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">  2; </w:t></w:r></w:p><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>greenLED</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">  3; </w:t></w:r></w:p><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>blinkTimes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">  5; </w:t></w:r></w:p><w:p><w:r><w:t>void setup() {</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">  </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>pinMode</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">(LED_BUILTIN, OUTPUT);         </w:t></w:r></w:p><w:p><w:r><w:t>}</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">void loop() { </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">   // Blink the Red Led 1 time for 1 second</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">   blink(</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>redLED</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">, 1); </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">   delay(1000); </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">   // Blink the Green Led 1 time for 1 second</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">   blink(</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>greenLED</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">, 1); </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">   delay(1000);</w:t></w:r></w:p><w:p><w:r><w:t>}</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\"> </w:t></w:r></w:p><w:p><w:r><w:t>void blink(</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>ledColor</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">, </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>blinkTimes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">) </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">{ </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">  if ((</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>blinkTimes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> &lt;=5) || (</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>blinkTimes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> &gt;5))</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">      for (</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> = 0; </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> &lt; </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>blinkTimes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">; </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">++) {                               </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">            </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>digitalWrite</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>ledColor</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">, HIGH);      </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">            delay(500);                      </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">            </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>digitalWrite</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>ledColor</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">, LOW);      </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">            delay(500);                    </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">      } </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">      </w:t></w:r></w:p><w:p><w:r><w:t>}</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>`;\n\nbody.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The edit rewrites the whole body: the 2-paragraph \"void setup(){...} void\n# loop(){...}\" sketch becomes a longer sketch that declares redLED/greenLED/\n# blinkTimes, calls a new blink() helper from loop(), and defines blink()\n# itself. Because so many runs/paragraphs are inserted (and the diff shows\n# exact <w:proofErr> spell-check bracketing around every camelCase\n# identifier), the most faithful way to reproduce the target OOXML is to\n# build that OOXML and hand it to Range.InsertXML, replacing the document's\n# whole story (Content) in one shot.\n\n$d = $word.ActiveDocument\n$null = $d.Paragraphs.Count\n\n$xml = @'\n<?xml version=\"1.0\" standalone=\"yes\"?><?mso-application progid=\"Word.Document\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>redLED</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">  2; </w:t></w:r></w:p><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>greenLED</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">  3; </w:t></w:r></w:p><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>blinkTimes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">  5; </w:t></w:r></w:p><w:p><w:r><w:t>void setup() {</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">  </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>pinMode</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">(LED_BUILTIN, OUTPUT);         </w:t></w:r></w:p><w:p><w:r><w:t>}</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">void loop() { </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">   // Blink the Red Led 1 time for 1 second</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">   blink(</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>redLED</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">, 1); </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">   delay(1000); </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">   // Blink the Green Led 1 time for 1 second</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">   blink(</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>greenLED</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">, 1); </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">   delay(1000);</w:t></w:r></w:p><w:p><w:r><w:t>}</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\"> </w:t></w:r></w:p><w:p><w:r><w:t>void blink(</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>ledColor</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">, </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>blinkTimes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">) </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">{ </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">  if ((</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>blinkTimes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> &lt;=5) || (</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>blinkTimes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> &gt;5))</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">      for (</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> = 0; </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> &lt; </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>blinkTimes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">; </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">++) {                               </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">            </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>digitalWrite</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>ledColor</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">, HIGH);      </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">            delay(500);                      </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">            </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>digitalWrite</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>ledColor</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">, LOW);      </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">            delay(500);                    </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">      } </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">      </w:t></w:r></w:p><w:p><w:r><w:t>}</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\n\n[void]$d.Content.InsertXML($xml)\n"}
